# "Taking back Chinna's changes"
#
# 1. Runmode column (C2:C5) on "Test Cases" reverts from "Y" back to "N".
# 2. Two test-case rows that had been removed are restored at the bottom
#    of the sheet (rows 6 & 7): ApplicationLinksValidationTest and
#    AppHeaderFooterLinkValidationTest.
# 3. Selection moves on to the newly-restored last row (C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Runmode column reverts to "N" for the existing rows ---
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"

# --- 2. Restore formatting (border/no-fill) for the two new rows by
#        copying it from the row immediately above, then overwrite with
#        the real values/heights/wrap for each new row. ---
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A5:D5").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6: ApplicationLinksValidationTest
$ws.Range("A6").Value = "ApplicationLinksValidationTest"
$ws.Range("B6").Value = "Validate below Application links " + [char]10 + "1. Web of Science" + [char]10 + "2.End Note" + [char]10 + "3.InCities" + [char]10 + "4.ScholarOne Abstracts" + [char]10 + "5.ScholarOne Manuscripts"
$ws.Range("C6").Value = "N"
$ws.Range("D6").Value = "SKIP"
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 90

# Row 7: AppHeaderFooterLinkValidationTest
$ws.Range("A7").Value = "AppHeaderFooterLinkValidationTest"
$ws.Range("B7").Value = "Validate Project Neon Header and Fooler links" + [char]10 + "1.Help" + [char]10 + "2.Cookie Policy" + [char]10 + "3.Privacy Statement" + [char]10 + "4.Terms of Use"
$ws.Range("C7").Value = "Y"
$ws.Range("D7").Value = "PASS"
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 75

# --- 3. Match the author's final selection ---
[void]$ws.Range("C7").Select()
